$wb = $excel.ActiveWorkbook

$renames = @{
    "train.data"      = "train_data"
    "test.data"       = "test_data"
    "train.data.prep" = "train_data_prep"
    "test.data.prep"  = "test_data_prep"
    "train.data.bc"   = "train_data_bc"
    "test.data.bc"    = "test_data_bc"
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($renames.ContainsKey($name)) {
        $ws.Name = $renames[$name]
    }
}
